# Bulgaria First League - odds update (21-04-2024 14:32)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Prepare the four new rows (249-252) so they inherit the same cell
# styles used by the existing data rows: bold/bordered/centered index
# column (A) and the date-formatted column (E). We copy the formats
# from row 248 (a representative existing data row) before writing
# values into the new rows.
# ---------------------------------------------------------------------
$ws.Range("A248").Copy() | Out-Null
$ws.Range("A249").PasteSpecial(-4122) | Out-Null
$ws.Range("A250").PasteSpecial(-4122) | Out-Null
$ws.Range("A251").PasteSpecial(-4122) | Out-Null
$ws.Range("A252").PasteSpecial(-4122) | Out-Null

$ws.Range("E248").Copy() | Out-Null
$ws.Range("E249").PasteSpecial(-4122) | Out-Null
$ws.Range("E250").PasteSpecial(-4122) | Out-Null
$ws.Range("E251").PasteSpecial(-4122) | Out-Null
$ws.Range("E252").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Row 248 — existing match record is updated with the finished result
# and closing odds (it also moves from kickoff 21/04 14:15 to 20/04 09:15).
# ---------------------------------------------------------------------
$ws.Range("B248").Value = 6978461
$ws.Range("E248").Value = 45402.38541666666
$ws.Range("F248").Value = "Pirin Blagoevgrad"
$ws.Range("G248").Value = "Cherno More Varna"
$ws.Range("H248").Value = 0
$ws.Range("I248").Value = 0
$ws.Range("J248").Value = "D"
$ws.Range("K248").Value = 5.75
$ws.Range("L248").Value = 3.75
$ws.Range("M248").Value = 1.6
$ws.Range("N248").Value = 6
$ws.Range("O248").Value = 3.6
$ws.Range("P248").Value = 1.6
$ws.Range("Q248").Value = 0.75
$ws.Range("R248").Value = 2.05
$ws.Range("S248").Value = 1.8
$ws.Range("T248").Value = 2.25
$ws.Range("U248").Value = 1.975
$ws.Range("V248").Value = 1.875
$ws.Range("W248").Value = -1
$ws.Range("X248").Value = 2.6
$ws.Range("Y248").Value = -1
$ws.Range("Z248").Value = 1.05
$ws.Range("AA248").Value = -1
$ws.Range("AB248").Value = -1
$ws.Range("AC248").Value = 0.875

# ---------------------------------------------------------------------
# Row 249 (new) — FC Hebar Pazardzhik vs Ludogorets Razgrad
# ---------------------------------------------------------------------
$ws.Range("A249").Value = 247
$ws.Range("B249").Value = 6978391
$ws.Range("C249").Value = "Bulgaria First League"
$ws.Range("D249").Value = "Bulgaria First League"
$ws.Range("E249").Value = 45402.48958333334
$ws.Range("F249").Value = "FC Hebar Pazardzhik"
$ws.Range("G249").Value = "Ludogorets Razgrad"
$ws.Range("H249").Value = 0
$ws.Range("I249").Value = 3
$ws.Range("J249").Value = "A"
$ws.Range("K249").Value = 10
$ws.Range("L249").Value = 6
$ws.Range("M249").Value = 1.25
$ws.Range("N249").Value = 12
$ws.Range("O249").Value = 7
$ws.Range("P249").Value = 1.222
$ws.Range("Q249").Value = 1.75
$ws.Range("R249").Value = 2.05
$ws.Range("S249").Value = 1.8
$ws.Range("T249").Value = 3
$ws.Range("U249").Value = 1.9
$ws.Range("V249").Value = 1.95
$ws.Range("W249").Value = -1
$ws.Range("X249").Value = -1
$ws.Range("Y249").Value = 0.222
$ws.Range("Z249").Value = -1
$ws.Range("AA249").Value = 0.8
$ws.Range("AB249").Value = 0
$ws.Range("AC249").Value = -0

# ---------------------------------------------------------------------
# Row 250 (new) — Levski Sofia vs Beroe
# ---------------------------------------------------------------------
$ws.Range("A250").Value = 248
$ws.Range("B250").Value = 6978462
$ws.Range("C250").Value = "Bulgaria First League"
$ws.Range("D250").Value = "Bulgaria First League"
$ws.Range("E250").Value = 45402.59375
$ws.Range("F250").Value = "Levski Sofia"
$ws.Range("G250").Value = "Beroe"
$ws.Range("H250").Value = 1
$ws.Range("I250").Value = 0
$ws.Range("J250").Value = "H"
$ws.Range("K250").Value = 1.444
$ws.Range("L250").Value = 4.2
$ws.Range("M250").Value = 7.5
$ws.Range("N250").Value = 1.285
$ws.Range("O250").Value = 5.25
$ws.Range("P250").Value = 9.5
$ws.Range("Q250").Value = -1.5
$ws.Range("R250").Value = 1.925
$ws.Range("S250").Value = 1.925
$ws.Range("T250").Value = 2.5
$ws.Range("U250").Value = 1.975
$ws.Range("V250").Value = 1.875
$ws.Range("W250").Value = 0.2849999999999999
$ws.Range("X250").Value = -1
$ws.Range("Y250").Value = -1
$ws.Range("Z250").Value = -1
$ws.Range("AA250").Value = 0.925
$ws.Range("AB250").Value = -1
$ws.Range("AC250").Value = 0.875

# ---------------------------------------------------------------------
# Row 251 (new) — Lokomotiv Plovdiv vs CSKA 1948 Sofia
# ---------------------------------------------------------------------
$ws.Range("A251").Value = 249
$ws.Range("B251").Value = 6978463
$ws.Range("C251").Value = "Bulgaria First League"
$ws.Range("D251").Value = "Bulgaria First League"
$ws.Range("E251").Value = 45403.38541666666
$ws.Range("F251").Value = "Lokomotiv Plovdiv"
$ws.Range("G251").Value = "CSKA 1948 Sofia"
$ws.Range("H251").Value = 1
$ws.Range("I251").Value = 1
$ws.Range("J251").Value = "D"
$ws.Range("K251").Value = 1.909
$ws.Range("L251").Value = 3.4
$ws.Range("M251").Value = 4
$ws.Range("N251").Value = 1.85
$ws.Range("O251").Value = 3.5
$ws.Range("P251").Value = 4.5
$ws.Range("Q251").Value = -0.5
$ws.Range("R251").Value = 1.85
$ws.Range("S251").Value = 2
$ws.Range("T251").Value = 2.5
$ws.Range("U251").Value = 1.975
$ws.Range("V251").Value = 1.875
$ws.Range("W251").Value = -1
$ws.Range("X251").Value = 2.5
$ws.Range("Y251").Value = -1
$ws.Range("Z251").Value = -1
$ws.Range("AA251").Value = 1
$ws.Range("AB251").Value = -1
$ws.Range("AC251").Value = 0.875

# ---------------------------------------------------------------------
# Row 252 (new) — Arda Kardzhali vs CSKA Sofia (fixture not yet played,
# so FTHG/FTAG/FTR are left blank, same as row 248 originally was).
# ---------------------------------------------------------------------
$ws.Range("A252").Value = 250
$ws.Range("B252").Value = 6978460
$ws.Range("C252").Value = "Bulgaria First League"
$ws.Range("D252").Value = "Bulgaria First League"
$ws.Range("E252").Value = 45403.59375
$ws.Range("F252").Value = "Arda Kardzhali"
$ws.Range("G252").Value = "CSKA Sofia"
$ws.Range("K252").Value = 5.5
$ws.Range("L252").Value = 3.6
$ws.Range("M252").Value = 1.65
$ws.Range("N252").Value = 5.25
$ws.Range("O252").Value = 3.6
$ws.Range("P252").Value = 1.666
$ws.Range("Q252").Value = 0.75
$ws.Range("R252").Value = 1.975
$ws.Range("S252").Value = 1.875
$ws.Range("T252").Value = 2.5
$ws.Range("U252").Value = 1.925
$ws.Range("V252").Value = 1.925
$ws.Range("W252").Value = 0
$ws.Range("X252").Value = 0
$ws.Range("Y252").Value = 0
$ws.Range("Z252").Value = 0
$ws.Range("AA252").Value = 0
